# ============================================================================
# edit.ps1 — Apply the "New crime data collected" weekly update to the
# 112th Precinct CompStat workbook (cs-en-us-112pct.xlsx).
#
# Changes:
#   1. Header text: Volume/Number bumped 43 -> 44; report week rolled
#      forward one week (10/24/2022-10/30/2022 -> 10/31/2022-11/6/2022).
#   2. Crime-complaint grid (rows 15-27, cols C:N): refreshed counts and
#      recomputed percent-change figures for the new week. Some cells flip
#      between a literal number and the sheet's textual "N/A" markers
#      ("0" for an undefined count, "***.*" for an undefined % change) —
#      those are written as text so they keep rendering the marker glyphs
#      instead of a numeric 0.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write one of the sheet's text "N/A" markers ("0" / "***.*") into a
# cell that currently holds a plain number. Forcing a text number format
# first keeps PowerShell/COM from re-coercing the literal "0" back into a
# numeric zero; flipping back to General afterwards keeps the visible
# formatting consistent with the rest of the table (text ignores
# NumberFormat for display purposes, so this is purely cosmetic/for
# consistency with neighboring text-marker cells).
function Set-TextMarker($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).NumberFormat = "General"
}

# Helper: the inverse — a cell currently holding a text marker now has a
# real count again, so give it back the grid's usual "#,##0" numeric format.
function Set-NumericFromText($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "#,##0"
    $ws.Range($addr).Value = $val
}

# --- Title text updates (Volume number + report week date range) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Crime-stat numeric cell updates (same type before/after: plain numeric overwrite) ---
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -53.846153846153
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = 65.217391304347
$ws.Range("L16").Value = 55.102040816326
$ws.Range("M16").Value = -26.213592233009
$ws.Range("N16").Value = -86.15664845173
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 71
$ws.Range("J17").Value = 56
$ws.Range("K17").Value = 26.785714285714
$ws.Range("L17").Value = 69.047619047619
$ws.Range("M17").Value = 44.897959183673
$ws.Range("N17").Value = -36.036036036036
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = 22.368421052631
$ws.Range("L18").Value = -7.920792079207
$ws.Range("M18").Value = -15.454545454545
$ws.Range("N18").Value = -92.301324503311
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -41.666666666666
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -3.225806451612
$ws.Range("I19").Value = 427
$ws.Range("J19").Value = 301
$ws.Range("K19").Value = 41.860465116279
$ws.Range("L19").Value = 69.444444444444
$ws.Range("M19").Value = 29.003021148036
$ws.Range("N19").Value = -50.40650406504
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 85
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = 44.067796610169
$ws.Range("L20").Value = 63.461538461538
$ws.Range("M20").Value = -1.162790697674
$ws.Range("N20").Value = -97.093023255813
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -56.521739130434
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = -13.636363636363
$ws.Range("I21").Value = 764
$ws.Range("J21").Value = 546
$ws.Range("K21").Value = 39.926739926739
$ws.Range("L21").Value = 51.888667992047
$ws.Range("M21").Value = 11.859443631039
$ws.Range("N21").Value = -86.513680494263
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 30
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 114.285714285714
$ws.Range("M22").Value = 50
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -12.121212121212
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = -22.222222222222
$ws.Range("I24").Value = 1498
$ws.Range("J24").Value = 1112
$ws.Range("K24").Value = 34.712230215827
$ws.Range("L24").Value = 67.749160134378
$ws.Range("M24").Value = 78.545887961859
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 9.090909090909
$ws.Range("I25").Value = 161
$ws.Range("J25").Value = 153
$ws.Range("K25").Value = 5.228758169934
$ws.Range("L25").Value = 38.793103448275
$ws.Range("M25").Value = -9.550561797752
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -20
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -50

# --- Cells changing FROM text marker back TO numeric ---
Set-NumericFromText $ws "C22" 1

# --- Cells changing FROM numeric TO text marker ("0" or "***.*") ---
Set-TextMarker $ws "C15" "0"
Set-TextMarker $ws "G15" "0"
Set-TextMarker $ws "H15" "***.*"
Set-TextMarker $ws "C16" "0"
Set-TextMarker $ws "C18" "0"
Set-TextMarker $ws "C20" "0"
Set-TextMarker $ws "G23" "0"
Set-TextMarker $ws "H23" "***.*"
Set-TextMarker $ws "D26" "0"
Set-TextMarker $ws "E26" "***.*"
Set-TextMarker $ws "D27" "0"
Set-TextMarker $ws "E27" "***.*"
